$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Probabilities")

# Row 2
$ws.Range("B2").Value = '2025-11-25T19:00:00'
$ws.Range("C2").Value = 'Нефтехимик'
$ws.Range("D2").Value = 'Авангард'
$ws.Range("E2").Value = 897814
$ws.Range("F2").Value = 'https://text.khl.ru/text/897814.html'
$ws.Range("G2").Value = 1.138889
$ws.Range("H2").Value = 5.538462
$ws.Range("I2").Value = 3.91689
$ws.Range("J2").Value = 1.307692
$ws.Range("K2").Value = 1.223291
$ws.Range("L2").Value = 4.727676
$ws.Range("M2").Value = 6.67735
$ws.Range("N2").Value = 23.62224
$ws.Range("O2").Value = 40.878477
$ws.Range("P2").Value = 64.50071699999999
$ws.Range("Q2").Value = -0.2
$ws.Range("R2").Value = 0.2
$ws.Range("S2").Value = 0.040577
$ws.Range("T2").Value = 0.059853
$ws.Range("U2").Value = 0.890176
$ws.Range("V2").Value = 0.155633
$ws.Range("W2").Value = 0.834974
$ws.Range("X2").Value = 0.291673
$ws.Range("Y2").Value = 0.6989340000000001
$ws.Range("Z2").Value = 0.453587
$ws.Range("AA2").Value = 0.5370200000000001
$ws.Range("AB2").Value = 0.614178
$ws.Range("AC2").Value = 0.376429
$ws.Range("AD2").Value = 0.750703
$ws.Range("AE2").Value = 0.239904
$ws.Range("AF2").Value = 0.345774
$ws.Range("AG2").Value = 0.654226
$ws.Range("AH2").Value = 0.125602
$ws.Range("AI2").Value = 0.874398
$ws.Range("AJ2").Value = 0.949327
$ws.Range("AK2").Value = 0.050673
$ws.Range("AL2").Value = 0.850458
$ws.Range("AM2").Value = 0.149542
$ws.Range("AN2").Value = 0.205009
$ws.Range("AO2").Value = 0.977089

# Row 3
$ws.Range("B3").Value = '2025-11-25T19:30:00'
$ws.Range("C3").Value = 'Спартак'
$ws.Range("D3").Value = 'Ак Барс'
$ws.Range("E3").Value = 897812
$ws.Range("F3").Value = 'https://text.khl.ru/text/897812.html'
$ws.Range("G3").Value = 4.588235
$ws.Range("H3").Value = 3.665934
$ws.Range("I3").Value = 4.852941
$ws.Range("J3").Value = 2.788167
$ws.Range("K3").Value = 3.688201
$ws.Range("L3").Value = 4.259438
$ws.Range("M3").Value = 8.254168999999999
$ws.Range("N3").Value = 36.153363
$ws.Range("O3").Value = 34.978634
$ws.Range("P3").Value = 71.131997
$ws.Range("Q3").Value = 0.2
$ws.Range("R3").Value = 0.113231
$ws.Range("S3").Value = 0.348824
$ws.Range("T3").Value = 0.141178
$ws.Range("U3").Value = 0.503965
$ws.Range("V3").Value = 0.043904
$ws.Range("W3").Value = 0.950063
$ws.Range("X3").Value = 0.10267
$ws.Range("Y3").Value = 0.891297
$ws.Range("Z3").Value = 0.19608
$ws.Range("AA3").Value = 0.797887
$ws.Range("AB3").Value = 0.319811
$ws.Range("AC3").Value = 0.674156
$ws.Range("AD3").Value = 0.460293
$ws.Range("AE3").Value = 0.533674
$ws.Range("AF3").Value = 0.882715
$ws.Range("AG3").Value = 0.117285
$ws.Range("AH3").Value = 0.712565
$ws.Range("AI3").Value = 0.287435
$ws.Range("AJ3").Value = 0.925683
$ws.Range("AK3").Value = 0.07431699999999999
$ws.Range("AL3").Value = 0.797502
$ws.Range("AM3").Value = 0.202498
$ws.Range("AN3").Value = 0.631793
$ws.Range("AO3").Value = 0.767918

# Row 4
$ws.Range("B4").Value = '2025-11-25T19:30:00'
$ws.Range("C4").Value = 'ХК Сочи'
$ws.Range("D4").Value = 'Лада'
$ws.Range("E4").Value = 897813
$ws.Range("F4").Value = 'https://text.khl.ru/text/897813.html'
$ws.Range("G4").Value = 1.166667
$ws.Range("H4").Value = 1.03125
$ws.Range("I4").Value = 1.233333
$ws.Range("J4").Value = 3.874015
$ws.Range("K4").Value = 2.520341
$ws.Range("L4").Value = 1.132292
$ws.Range("M4").Value = 2.197917
$ws.Range("N4").Value = 25.593223
$ws.Range("O4").Value = 23.641989
$ws.Range("P4").Value = 49.235212
$ws.Range("Q4").Value = -0.2
$ws.Range("R4").Value = -0.2
$ws.Range("S4").Value = 0.675169
$ws.Range("T4").Value = 0.172768
$ws.Range("U4").Value = 0.151996
$ws.Range("V4").Value = 0.5040829999999999
$ws.Range("W4").Value = 0.49585
$ws.Range("X4").Value = 0.696346
$ws.Range("Y4").Value = 0.303588
$ws.Range("Z4").Value = 0.836799
$ws.Range("AA4").Value = 0.163135
$ws.Range("AB4").Value = 0.922303
$ws.Range("AC4").Value = 0.07763100000000001
$ws.Range("AD4").Value = 0.966919
$ws.Range("AE4").Value = 0.033015
$ws.Range("AF4").Value = 0.716851
$ws.Range("AG4").Value = 0.283149
$ws.Range("AH4").Value = 0.461394
$ws.Range("AI4").Value = 0.538606
$ws.Range("AJ4").Value = 0.312776
$ws.Range("AK4").Value = 0.6872239999999999
$ws.Range("AL4").Value = 0.106172
$ws.Range("AM4").Value = 0.893828
$ws.Range("AN4").Value = 0.944574
$ws.Range("AO4").Value = 0.539866
